$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price-column values are purely numeric-looking strings (e.g. "0.999").
# The source workbook stores every data cell as literal text (t="inlineStr"),
# so for those cells we briefly force a text number format before assigning the
# value (otherwise Excel auto-converts the input into a real number), then
# restore the cell's original style so formatting is unaffected.

# Row 2
$ws.Range("D2").Value = '38.780.09'
$ws.Range("E2").Value = '  +1.90%  '

# Row 3
$ws.Range("D3").Value = '2.094.28'
$ws.Range("E3").Value = '  +0.17%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.57'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.16%  '

# Row 6
$ws.Range("E6").Value = '  +0.41%  '

# Row 7
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.14'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +0.41%  '

# Row 8
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.387'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +1.81%  '

# Row 10
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0845'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +0.38%  '

# Row 11
$ws.Range("E11").Value = '  -0.58%  '

# Row 12
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.36'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +5.09%  '

# Row 13
$ws.Range("D13").Value = '2.409.07'
$ws.Range("E13").Value = '  +0.39%  '

# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.08'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -0.58%  '

# Row 15
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.807'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +4.21%  '

# Row 16
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").Value = '2.092.57'
$ws.Range("E17").Value = '  +0.45%  '

# Row 18
$ws.Range("D18").Value = '38.691.69'
$ws.Range("E18").Value = '  +1.82%  '

# Row 19
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.04'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +2.84%  '

# Row 20
$ws.Range("E20").Value = '  +0.74%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  +0.40%  '

# Row 22
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.79'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +1.71%  '

# Row 23
$ws.Range("E23").Value = '  -0.43%  '

# Row 24
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -2.27%  '

# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.89%  '

# Row 26
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.33'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +0.87%  '

# Row 27
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.55'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +1.76%  '

# Row 28
$ws.Range("E28").Value = '  +5.46%  '

# Row 29
$ws.Range("E29").Value = '  +6.73%  '

# Row 30
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.33'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +2.00%  '

# Row 31
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.46'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +2.34%  '

# Row 32
$ws.Range("E32").Value = '  +0.69%  '

# Row 33
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.52'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +2.09%  '

# Row 34
$ws.Range("E34").Value = '  +0.84%  '

# Row 35
$ws.Range("E35").Value = '  +0.80%  '

# Row 36
$ws.Range("E36").Value = '  +2.03%  '

# Row 37
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("E38").Value = '  +1.54%  '

# Row 39
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -0.13%  '

# Row 40
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.19'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +0.34%  '

# Row 41
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0230'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +5.22%  '

# Row 42
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.08'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +1.20%  '

# Row 43
$ws.Range("D43").Value = '1.533.37'
$ws.Range("E43").Value = '  -0.90%  '

# Row 45
$ws.Range("E45").Value = '  +1.38%  '

# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.14'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +2.07%  '

# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.69'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +6.35%  '

# Row 48
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.13'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -0.74%  '

# Row 49
$ws.Range("E49").Value = '  +1.51%  '

# Row 50
$ws.Range("E50").Value = '  -0.48%  '

# Row 51
$ws.Range("D51").Value = '2.290.14'
$ws.Range("E51").Value = '  +0.10%  '
